$d = $word.ActiveDocument

# Locate the paragraph that ends with the known tail text
# ("...Configured authorization in the startup.cs file.").
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Configured authorization in the startup\.cs file\.") {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find target paragraph"
}

$target = $d.Paragraphs.Item($targetIndex)

# Collapse a range to the very end of the target paragraph (right before its
# paragraph mark). Inserting a new paragraph here creates a sibling list item
# that inherits the ListParagraph style + numbering (numId 1) from $target.
$ip = $target.Range
$ip.Collapse(0)  # wdCollapseEnd -> 0
$ip.InsertParagraphAfter()

# --- New paragraph 1: "Made the shopping cart navigation bar item clickable." ---
$para1 = $d.Paragraphs.Item($targetIndex + 1)
$r1 = $d.Range($para1.Range.Start, $para1.Range.Start)
$r1.InsertAfter("Made the shopping cart navigation bar item clickable.")
$r1.Collapse(0)
$r1.InsertParagraphAfter()

# --- New paragraph 2: "Added a button to clear the shopping cart." (3 runs) ---
$para2 = $d.Paragraphs.Item($targetIndex + 2)
$r2 = $d.Range($para2.Range.Start, $para2.Range.Start)
$r2.InsertAfter("Added")
$r2.Collapse(0)
$r2.InsertAfter(" a")
$r2.Collapse(0)
$r2.InsertAfter(" button to clear the shopping cart.")
$r2.Collapse(0)
$r2.InsertParagraphAfter()

# --- New paragraph 3: "Creating a link to decrease the amount in the shopping cart." ---
$para3 = $d.Paragraphs.Item($targetIndex + 3)
$r3 = $d.Range($para3.Range.Start, $para3.Range.Start)
$r3.InsertAfter("Creating a link to decrease the amount in the shopping cart.")
